$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Repayment Schedule")

# Insert a new blank column before column N (14) - shifts old N,O,P to O,P,Q
$ws.Columns.Item(14).Insert()

# New column N width (displayed as 10 characters; ColumnWidth input needs the
# 5/6-character offset Excel applies internally so the stored width is 10)
$ws.Columns.Item(14).ColumnWidth = 9.166666666666666

# Activate this sheet / select the new active cell to match target view state
$ws.Activate()
$ws.Range("M14").Select()
